$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 350
$ws.Range("I2").Value = 350
$ws.Range("K2").Value = 350
$ws.Range("M2").Value = -237

$ws.Range("H28").Value = 721.9167
$ws.Range("I28").Value = 574.75
$ws.Range("K28").Value = 574.75
$ws.Range("M28").Value = -89.75

$ws.Range("H33").Value = 846.8421
$ws.Range("I33").Value = 960.9231
$ws.Range("K33").Value = 960.9231
$ws.Range("M33").Value = -731.9231

$ws.Range("H86").Value = 154324690
$ws.Range("I86").Value = 266668530
$ws.Range("K86").Value = 266668530
$ws.Range("M86").Value = -266667407

$ws.Range("H89").Value = 154324690
$ws.Range("I89").Value = 266668530
$ws.Range("K89").Value = 1333342650
$ws.Range("M89").Value = -1333337034

$ws.Range("H116").Value = 17864342
$ws.Range("I116").Value = 35717900
$ws.Range("K116").Value = 35717900
$ws.Range("M116").Value = -35714458

$ws.Range("H121").Value = 5692
$ws.Range("J121").Value = 5692
$ws.Range("L121").Value = 17076
$ws.Range("N121").Value = -20570

$ws.Range("H131").Value = 2470.4443

$ws.Range("H132").Value = 2288.6287
$ws.Range("I132").Value = 1724.5518
$ws.Range("J132").Value = 5015
$ws.Range("K132").Value = 5173.6554
$ws.Range("L132").Value = 15045
$ws.Range("M132").Value = -2643.6554
$ws.Range("N132").Value = -20105

$ws.Range("H135").Value = 435336.8
$ws.Range("I135").Value = 526807.8
$ws.Range("K135").Value = 4741270.2
$ws.Range("M135").Value = -4738735.2

$ws.Range("H137").Value = 3230.077
$ws.Range("I137").Value = 5079.2
$ws.Range("J137").Value = 2074.375
$ws.Range("K137").Value = 15237.6
$ws.Range("L137").Value = 6223.125
$ws.Range("M137").Value = -12687.6
$ws.Range("N137").Value = -11323.125

$ws.Range("H138").Value = 8658
$ws.Range("I138").Value = 1409.6666
$ws.Range("J138").Value = 16812.375
$ws.Range("K138").Value = 4228.9998
$ws.Range("L138").Value = 50437.125
$ws.Range("M138").Value = 911.0002000000004
$ws.Range("N138").Value = -60717.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1472721.6
$ws.Range("I32").Value = 1507997.9
$ws.Range("K32").Value = 1507997.9
$ws.Range("M32").Value = -1507710.9

$ws.Range("H45").Value = 3711.6875
$ws.Range("I45").Value = 1494.25
$ws.Range("J45").Value = 5929.125
$ws.Range("K45").Value = 1494.25
$ws.Range("L45").Value = 5929.125
$ws.Range("M45").Value = -1117.25
$ws.Range("N45").Value = -6683.125

$ws.Range("H61").Value = 6173.638
$ws.Range("I61").Value = 2818.0293
$ws.Range("K61").Value = 2818.0293
$ws.Range("M61").Value = -2606.0293

$ws.Range("H69").Value = 250000
$ws.Range("J69").Value = 250000
$ws.Range("L69").Value = 250000
$ws.Range("N69").Value = -251498

$ws.Range("H72").Value = 250000
$ws.Range("J72").Value = 250000
$ws.Range("L72").Value = 750000
$ws.Range("N72").Value = -757488

$ws.Range("H74").Value = 103531.625
$ws.Range("I74").Value = 201882.5
$ws.Range("K74").Value = 201882.5
$ws.Range("M74").Value = -201008.5

$ws.Range("H77").Value = 103531.625
$ws.Range("I77").Value = 201882.5
$ws.Range("K77").Value = 1009412.5
$ws.Range("M77").Value = -1005044.5

$ws.Range("H97").Value = 3477319.2
$ws.Range("J97").Value = 10431061
$ws.Range("L97").Value = 10431061
$ws.Range("N97").Value = -10432053

$ws.Range("H102").Value = 1498.3334
$ws.Range("I102").Value = 1742
$ws.Range("J102").Value = 1193.75
$ws.Range("K102").Value = 1742
$ws.Range("L102").Value = 1193.75
$ws.Range("M102").Value = -120
$ws.Range("N102").Value = -4437.75

$ws.Range("H136").Value = 6173.638
$ws.Range("I136").Value = 2818.0293
$ws.Range("K136").Value = 8454.0879
$ws.Range("M136").Value = -5904.0879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 58749.6
$ws.Range("J110").Value = 58749.6
$ws.Range("L110").Value = 58749.6
$ws.Range("N110").Value = -66929.60000000001

$ws.Range("H134").Value = 7973.185
$ws.Range("I134").Value = 3094.75
$ws.Range("K134").Value = 9284.25
$ws.Range("M134").Value = -6749.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 31975
$ws.Range("J36").Value = 49950
$ws.Range("L36").Value = 49950
$ws.Range("N36").Value = -50726

$ws.Range("H40").Value = 31975
$ws.Range("J40").Value = 49950
$ws.Range("L40").Value = 49950
$ws.Range("N40").Value = -50270

$ws.Range("H58").Value = 17249340
$ws.Range("I58").Value = 55556572
$ws.Range("J58").Value = 11086.65
$ws.Range("K58").Value = 55556572
$ws.Range("L58").Value = 11086.65
$ws.Range("M58").Value = -55556369
$ws.Range("N58").Value = -11492.65

$ws.Range("H132").Value = 10062
$ws.Range("I132").Value = 3332.3333
$ws.Range("J132").Value = 11615
$ws.Range("K132").Value = 9996.999899999999
$ws.Range("L132").Value = 34845
$ws.Range("M132").Value = -7466.999899999999
$ws.Range("N132").Value = -39905

$ws.Range("H136").Value = 17249340
$ws.Range("I136").Value = 55556572
$ws.Range("J136").Value = 11086.65
$ws.Range("K136").Value = 166669716
$ws.Range("L136").Value = 33259.95
$ws.Range("M136").Value = -166667166
$ws.Range("N136").Value = -38359.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 126290.31
$ws.Range("I2").Value = 57.545456
$ws.Range("J2").Value = 404002.4
$ws.Range("K2").Value = 345.272736
$ws.Range("L2").Value = 2424014.4
$ws.Range("M2").Value = -232.272736
$ws.Range("N2").Value = -2424240.4

$ws.Range("H7").Value = 242.5
$ws.Range("I7").Value = 242.5
$ws.Range("K7").Value = 727.5
$ws.Range("M7").Value = -615.5

$ws.Range("H92").Value = 6994459.5
$ws.Range("J92").Value = 6994459.5
$ws.Range("L92").Value = 20983378.5
$ws.Range("N92").Value = -20985874.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8454.817999999999
$ws.Range("I7").Value = 6000
$ws.Range("J7").Value = 8700.299999999999
$ws.Range("K7").Value = 6000
$ws.Range("L7").Value = 8700.299999999999
$ws.Range("M7").Value = -5888
$ws.Range("N7").Value = -8924.299999999999

$ws.Range("H46").Value = 4148
$ws.Range("I46").Value = 4000
$ws.Range("K46").Value = 4000
$ws.Range("M46").Value = -3812

$ws.Range("H68").Value = 5672.1113
$ws.Range("J68").Value = 6175
$ws.Range("L68").Value = 6175
$ws.Range("N68").Value = -7673

$ws.Range("H71").Value = 5672.1113
$ws.Range("J71").Value = 6175
$ws.Range("L71").Value = 30875
$ws.Range("N71").Value = -38363

$ws.Range("H100").Value = 3210.12
$ws.Range("I100").Value = 3019
$ws.Range("K100").Value = 3019
$ws.Range("M100").Value = -2478

$ws.Range("H126").Value = 8454.817999999999
$ws.Range("I126").Value = 6000
$ws.Range("J126").Value = 8700.299999999999
$ws.Range("K126").Value = 18000
$ws.Range("L126").Value = 26100.9
$ws.Range("M126").Value = -15530
$ws.Range("N126").Value = -31040.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6670065.5
$ws.Range("I81").Value = 1401.7858
$ws.Range("J81").Value = 12505146
$ws.Range("K81").Value = 2803.5716
$ws.Range("L81").Value = 25010292
$ws.Range("M81").Value = -1742.5716
$ws.Range("N81").Value = -25012414

$ws.Range("H84").Value = 6670065.5
$ws.Range("I84").Value = 1401.7858
$ws.Range("J84").Value = 12505146
$ws.Range("K84").Value = 14017.858
$ws.Range("L84").Value = 125051460
$ws.Range("M84").Value = -8713.858
$ws.Range("N84").Value = -125062068

$ws.Range("H132").Value = 22733882
$ws.Range("I132").Value = 41675960
$ws.Range("J132").Value = 3389.5
$ws.Range("K132").Value = 125027880
$ws.Range("L132").Value = 10168.5
$ws.Range("M132").Value = -125025350
$ws.Range("N132").Value = -15228.5

$ws.Range("H136").Value = 35719416
$ws.Range("I136").Value = 83333970
$ws.Range("J136").Value = 8500.625
$ws.Range("K136").Value = 250001910
$ws.Range("L136").Value = 25501.875
$ws.Range("M136").Value = -249999360
$ws.Range("N136").Value = -30601.875
